# 2040_MT.xlsx — "Updated Results with corrected code"
#
# Changes applied to Sheet1:
#   1. D3: the stray numeric 0 is cleared back to blank (no longer a
#      computed/placeholder value for Hydrogen's Non-metallic minerals cell).
#   2. A7: the row labelled "Other" is renamed to "Biogas" (a fuel/feedstock
#      that was missing from the breakdown).
#   3. A new row 8, labelled "Other", is appended below it - reusing the
#      same look (bold, centered, bordered label in column A) as the other
#      rows and the same blank/"0" pattern as the existing fuel rows
#      (Biomass/Other) in columns B-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) D3 goes back to blank.
$ws.Range("D3").ClearContents()

# 2) Row 7's label "Other" becomes "Biogas".
$ws.Range("A7").Value = "Biogas"

# 3) Insert the new "Other" row (row 8), cloning row 7's formatting so the
#    label keeps the bold/centered/bordered style, then filling in values:
#    A8 = "Other" (text), B8/C8 = blank, D8 = 0 (numeric), matching the
#    existing Biomass/Biogas rows.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 0
